$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87, shifting existing rows 87-180 down to 88-181.
$ws.Rows("87:87").Insert()

# Populate the newly inserted row 87 with the new record's data.
$ws.Range("A87").Value = 7
$ws.Range("B87").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C87").Value = "Ñuble"
$ws.Range("D87").Value = 44539
$ws.Range("E87").Value = 16
$ws.Range("F87").Value = 100112043
$ws.Range("G87").Value = "Pepino ensalada"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 240
$ws.Range("K87").Value = 6500
$ws.Range("L87").Value = 7000
$ws.Range("M87").Value = 6750
$ws.Range("N87").Value = "$/caja 80 unidades"
$ws.Range("O87").Value = "Región del Maule"
$ws.Range("P87").Value = 84
$ws.Range("Q87").Value = 80
$ws.Range("R87").Value = "Hortaliza"
